# sep 30th new update
# Add the Sep 30 (row 13) Covid stats: previous-cases delta (E), new total (F),
# and deaths (G). D13 ("#Previous cases" total, 1660) was already filled in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E13").Value = 137
$ws.Range("E13").Style = "Bad"
$ws.Range("F13").Value = 1797
$ws.Range("G13").Value = 6

# Move the selection to where the user's cursor landed after the entry.
$ws.Range("G14").Select() | Out-Null
